## contingencies with rene fine
## Insert two new rows ("line7", "line8") after the existing "line6" row
## (shifting the extr1..extr8 rows down by two), then update the
## from_bus/to_bus/in_service values for the two new rows plus several of
## the existing extr* rows, matching the target workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- insert two blank rows at 8:9, pushing extr1..extr8 down to 10..17 ---
$ws.Rows("8:9").Insert()

# Copy the formatting of the row right above (line6, row 7) down into the
# two freshly-inserted rows so the index column keeps its existing style.
$ws.Range("A7:E7").Copy()
$ws.Range("A8:E9").PasteSpecial(-4122)   # xlPasteFormats

# --- new row 8: line7 ---
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# --- new row 9: line8 ---
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# --- fix up the index column + in_service flags for the shifted extr rows ---
# row 10: extr1 (index + from/to unchanged, in_service False -> True)
$ws.Range("A10").Value = 8
$ws.Range("E10").Value = $true

# row 11: extr2 (index only)
$ws.Range("A11").Value = 9

# row 12: extr3 (index only, values unchanged)
$ws.Range("A12").Value = 10

# row 13: extr4 (index, in_service False -> True)
$ws.Range("A13").Value = 11
$ws.Range("E13").Value = $true

# row 14: extr5 (index, in_service False -> True)
$ws.Range("A14").Value = 12
$ws.Range("E14").Value = $true

# row 15: extr6 (index, in_service False -> True)
$ws.Range("A15").Value = 13
$ws.Range("E15").Value = $true

# row 16: extr7 (index only, stays False)
$ws.Range("A16").Value = 14

# row 17: extr8 (index, in_service True -> False)
$ws.Range("A17").Value = 15
$ws.Range("E17").Value = $false
